# Scheduled-runner update: refresh cached Market Board price / profit
# figures (currentAveragePrice*, LevePrice*, LeveProfit*) across several
# class sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR) of the Bahamut Profits
# workbook. Values below mirror newly-pulled Universalis data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 30228.5
$ws.Range("J3").Value = 30228.5
$ws.Range("L3").Value = 30228.5
$ws.Range("N3").Value = -30456.5

$ws.Range("H43").Value = 589030.3
$ws.Range("I43").Value = 783.8889
$ws.Range("J43").Value = 1250807.5
$ws.Range("K43").Value = 783.8889
$ws.Range("L43").Value = 1250807.5
$ws.Range("M43").Value = -714.8889
$ws.Range("N43").Value = -1250945.5

$ws.Range("H88").Value = 1444712.4
$ws.Range("I88").Value = 5251.5
$ws.Range("J88").Value = 1764592.5
$ws.Range("K88").Value = 5251.5
$ws.Range("L88").Value = 1764592.5
$ws.Range("M88").Value = -4845.5
$ws.Range("N88").Value = -1765404.5

$ws.Range("H91").Value = 1444712.4
$ws.Range("I91").Value = 5251.5
$ws.Range("J91").Value = 1764592.5
$ws.Range("K91").Value = 5251.5
$ws.Range("L91").Value = 1764592.5
$ws.Range("M91").Value = -3847.5
$ws.Range("N91").Value = -1767400.5

$ws.Range("H102").Value = 30228.5
$ws.Range("J102").Value = 30228.5
$ws.Range("L102").Value = 30228.5
$ws.Range("N102").Value = -36718.5

$ws.Range("H137").Value = 3824.25
$ws.Range("I137").Value = 3518.8
$ws.Range("J137").Value = 4333.3335
$ws.Range("K137").Value = 10556.4
$ws.Range("L137").Value = 13000.0005
$ws.Range("M137").Value = -8006.400000000001
$ws.Range("N137").Value = -18100.0005

$ws.Range("H140").Value = 47414.285
$ws.Range("J140").Value = 47414.285
$ws.Range("L140").Value = 47414.285
$ws.Range("N140").Value = -57774.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3287.4285
$ws.Range("I45").Value = 2622.4
$ws.Range("J45").Value = 4950
$ws.Range("K45").Value = 2622.4
$ws.Range("L45").Value = 4950
$ws.Range("M45").Value = -2245.4
$ws.Range("N45").Value = -5704

$ws.Range("H63").Value = 2270.3333
$ws.Range("I63").Value = 2452.5
$ws.Range("J63").Value = 1906
$ws.Range("K63").Value = 2452.5
$ws.Range("L63").Value = 1906
$ws.Range("M63").Value = -1766.5
$ws.Range("N63").Value = -3278

$ws.Range("H64").Value = 28600
$ws.Range("J64").Value = 28600
$ws.Range("L64").Value = 28600
$ws.Range("N64").Value = -29096

$ws.Range("H66").Value = 2270.3333
$ws.Range("I66").Value = 2452.5
$ws.Range("J66").Value = 1906
$ws.Range("K66").Value = 12262.5
$ws.Range("L66").Value = 9530
$ws.Range("M66").Value = -8830.5
$ws.Range("N66").Value = -16394

$ws.Range("H67").Value = 28600
$ws.Range("J67").Value = 28600
$ws.Range("L67").Value = 28600
$ws.Range("N67").Value = -30316

$ws.Range("H74").Value = 782.6429000000001
$ws.Range("I74").Value = 726.6957
$ws.Range("J74").Value = 1040
$ws.Range("K74").Value = 726.6957
$ws.Range("L74").Value = 1040
$ws.Range("M74").Value = 147.3043
$ws.Range("N74").Value = -2788

$ws.Range("H77").Value = 782.6429000000001
$ws.Range("I77").Value = 726.6957
$ws.Range("J77").Value = 1040
$ws.Range("K77").Value = 3633.4785
$ws.Range("L77").Value = 5200
$ws.Range("M77").Value = 734.5214999999998
$ws.Range("N77").Value = -13936

$ws.Range("H88").Value = 1937.6
$ws.Range("I88").Value = 1671.4166
$ws.Range("J88").Value = 3002.3333
$ws.Range("K88").Value = 1671.4166
$ws.Range("L88").Value = 3002.3333
$ws.Range("M88").Value = -1265.4166
$ws.Range("N88").Value = -3814.3333

$ws.Range("H91").Value = 1937.6
$ws.Range("I91").Value = 1671.4166
$ws.Range("J91").Value = 3002.3333
$ws.Range("K91").Value = 1671.4166
$ws.Range("L91").Value = 3002.3333
$ws.Range("M91").Value = -267.4166
$ws.Range("N91").Value = -5810.3333

$ws.Range("H122").Value = 1496.2858
$ws.Range("I122").Value = 1496.2858
$ws.Range("K122").Value = 4488.857400000001
$ws.Range("M122").Value = -2038.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3508.2273
$ws.Range("I31").Value = 1418.7333
$ws.Range("J31").Value = 7985.7144
$ws.Range("K31").Value = 1418.7333
$ws.Range("L31").Value = 7985.7144
$ws.Range("M31").Value = -1123.7333
$ws.Range("N31").Value = -8575.714400000001

$ws.Range("H34").Value = 3508.2273
$ws.Range("I34").Value = 1418.7333
$ws.Range("J34").Value = 7985.7144
$ws.Range("K34").Value = 1418.7333
$ws.Range("L34").Value = 7985.7144
$ws.Range("M34").Value = -1216.7333
$ws.Range("N34").Value = -8389.714400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1683.7142
$ws.Range("I5").Value = 892.9091
$ws.Range("J5").Value = 4583.3335
$ws.Range("K5").Value = 2678.7273
$ws.Range("L5").Value = 13750.0005
$ws.Range("M5").Value = -2566.7273
$ws.Range("N5").Value = -13974.0005

$ws.Range("H112").Value = 3017
$ws.Range("J112").Value = 3810
$ws.Range("L112").Value = 11430
$ws.Range("N112").Value = -13646

$ws.Range("H121").Value = 1464.2858
$ws.Range("I121").Value = 315
$ws.Range("J121").Value = 1655.8334
$ws.Range("K121").Value = 945
$ws.Range("L121").Value = 4967.5002
$ws.Range("M121").Value = 365
$ws.Range("N121").Value = -7587.5002

$ws.Range("H135").Value = 1683.7142
$ws.Range("I135").Value = 892.9091
$ws.Range("J135").Value = 4583.3335
$ws.Range("K135").Value = 8036.1819
$ws.Range("L135").Value = 41250.0015
$ws.Range("M135").Value = -5501.1819
$ws.Range("N135").Value = -46320.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 35260
$ws.Range("J95").Value = 35260
$ws.Range("L95").Value = 35260
$ws.Range("N95").Value = -40752

$ws.Range("H105").Value = 28000
$ws.Range("I105").Value = 28000
$ws.Range("K105").Value = 28000
$ws.Range("M105").Value = -24506

$ws.Range("H123").Value = 29236
$ws.Range("J123").Value = 29236
$ws.Range("L123").Value = 29236
$ws.Range("N123").Value = -34136

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 12831
$ws.Range("J101").Value = 12831
$ws.Range("L101").Value = 12831
$ws.Range("N101").Value = -19321

$ws.Range("H122").Value = 12208.75
$ws.Range("I122").Value = 13530.714
$ws.Range("J122").Value = 2955
$ws.Range("K122").Value = 40592.142
$ws.Range("L122").Value = 8865
$ws.Range("M122").Value = -38142.142
$ws.Range("N122").Value = -13765

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H122").Value = 1566.4445
$ws.Range("I122").Value = 1456.2222
$ws.Range("J122").Value = 1676.6666
$ws.Range("K122").Value = 4368.6666
$ws.Range("L122").Value = 5029.9998
$ws.Range("M122").Value = -1918.6666
$ws.Range("N122").Value = -9929.9998

$ws.Range("H123").Value = 38098.777
$ws.Range("J123").Value = 38098.777
$ws.Range("L123").Value = 38098.777
$ws.Range("N123").Value = -47898.777
